$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.193811058998108
$ws.Range("B1").Value = 1.44761860370636
$ws.Range("C1").Value = 6.771073818206787
$ws.Range("D1").Value = 2.095736742019653
$ws.Range("E1").Value = 0.9225821495056152
